$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.530.66'
$ws.Range("E2").Value = '  -7.34%  '
$ws.Range("D3").Value = '1.691.02'
$ws.Range("E3").Value = '  -5.86%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5130'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -13.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.006'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2670'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '22.09'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06322'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07379'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.25%  '
$ws.Range("D12").Value = '1.694.98'
$ws.Range("E12").Value = '  -5.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.533'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5794'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.06%  '
$ws.Range("D15").Value = '1.920.80'
$ws.Range("E15").Value = '  -5.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008608'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -13.48%  '
$ws.Range("D18").Value = '26.579.47'
$ws.Range("E18").Value = '  -7.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.994'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.007'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("E21").Value = '  -4.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '187.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.272'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.542'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1176'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.359'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05819'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.73%  '
$ws.Range("E31").Value = '  -6.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.532'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.526'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.658'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.82%  '
$ws.Range("E35").Value = '  -2.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.5991'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.97%  '
$ws.Range("E37").Value = '  -5.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.672'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01620'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.91%  '
$ws.Range("D40").Value = '1.097.75'
$ws.Range("E40").Value = '  -4.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8655'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.895'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.07%  '
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("D45").Value = '1.848.29'
$ws.Range("E45").Value = '  -5.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000114'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.057'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05239'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4320'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.53%  '
